# "Added normalization to WT": insert a new top row holding the
# normalization-to-wild-type (BY4743) values, pushing the existing
# strain rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1; every existing row shifts down.
$ws.Rows("1:1").Insert()

# Populate the new row 1 with the BY4743 (wild-type) normalization data.
$ws.Range("A1").Value = "BY4743"

$ws.Range("C1").Value = 0.08
$ws.Range("D1").Value = "(0.15)"
$ws.Range("D1").NumberFormat = "@"

$ws.Range("E1").Value = 0.05
$ws.Range("F1").Value = "(0.09)"
$ws.Range("F1").NumberFormat = "@"

$ws.Range("G1").Value = 0.06
$ws.Range("H1").Value = "(0.32)"
$ws.Range("H1").NumberFormat = "@"

$ws.Range("I1").Value = 0.03
$ws.Range("J1").Value = "(0.31)"
$ws.Range("J1").NumberFormat = "@"

$ws.Range("K1").Value = 0.22

# Leave the selection on B14, matching the saved worksheet view.
$ws.Range("B14").Select()
